$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Log")

# Enter the measured value for row 33 (column C), which was previously blank.
# This feeds the dependent formulas in D33/E33 and the AVERAGEIF results in
# column F for the "15.Mix.extended" group (rows 32-37).
$ws.Range("C33").Value = 0.71875

# Recalculate so the dependent formula cells (E33, F32:F37, etc.) update.
$excel.Calculate()

# Move the active selection to C34, matching where the user left off after
# entering the new data point.
$ws.Activate()
$ws.Range("C34").Select()
